# Apply Chinese-name translations (and a couple of garbled machine-translation
# artifacts that appear verbatim in the target workbook) to column B
# ("Translated") of Sheet1, rows 2-56. Column A ("Original") is left untouched.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Cells.Item(2, 2).Value = "安徽交通控股集团有限公司"
$ws.Cells.Item(3, 2).Value = "北京未来科学园发展集团有限公司"
$ws.Cells.Item(4, 2).Value = "北京五谷集团有限公司"
$ws.Cells.Item(5, 2).Value = "北京国有资产管理有限公司"
$ws.Cells.Item(6, 2).Value = "长春城市发展投资控股（集团）有限公司"
$ws.Cells.Item(7, 2).Value = "Cheng度airport xing成investment group co., Ltd."
$ws.Cells.Item(8, 2).Value = "Cheng度airport xing成investment group co., Ltd."
$ws.Cells.Item(9, 2).Value = "成都交通投资集团有限公司"
$ws.Cells.Item(10, 2).Value = "中国华融资产管理有限公司"
$ws.Cells.Item(11, 2).Value = "重庆市南岸城市建设发展（集团）有限公司"
$ws.Cells.Item(12, 2).Value = "DA连D ETA holding co., Ltd."
$ws.Cells.Item(13, 2).Value = "甘肃省公路航空旅游投资集团有限公司"
$ws.Cells.Item(14, 2).Value = "GU案G栋he ng见investment holding co., Ltd."
$ws.Cells.Item(15, 2).Value = "GU案G栋he ng见investment holding co., Ltd."
$ws.Cells.Item(16, 2).Value = "GU案G系L IU周dong成investment development group co., Ltd."
$ws.Cells.Item(17, 2).Value = "广州开发区金融控股集团有限公司"
$ws.Cells.Item(18, 2).Value = "H安检GG UO头group co., Ltd."
$ws.Cells.Item(19, 2).Value = "河北建设投资集团有限公司"
$ws.Cells.Item(20, 2).Value = "河北建设投资集团有限公司"
$ws.Cells.Item(21, 2).Value = "河北省国有资产控股经营有限公司"
$ws.Cells.Item(22, 2).Value = "合肥工业投资控股（集团）有限公司"
$ws.Cells.Item(23, 2).Value = "河南投资集团有限公司"
$ws.Cells.Item(24, 2).Value = "河南投资集团有限公司"
$ws.Cells.Item(25, 2).Value = "衡阳市建设投资有限公司"
$ws.Cells.Item(26, 2).Value = "淮安发展控股有限公司"
$ws.Cells.Item(27, 2).Value = "hu AI按traffic holding co., Ltd."
$ws.Cells.Item(28, 2).Value = "江苏汉瑞投资控股有限公司"
$ws.Cells.Item(29, 2).Value = "Jiang苏jin管investment and development group co., Ltd."
$ws.Cells.Item(30, 2).Value = "江西铁路投资集团公司"
$ws.Cells.Item(31, 2).Value = "Kunming Dian吃investment co., Ltd."
$ws.Cells.Item(32, 2).Value = "昆明交通投资有限公司"
$ws.Cells.Item(33, 2).Value = "昆明交通投资有限公司"
$ws.Cells.Item(34, 2).Value = "昆山交通发展控股集团有限公司"
$ws.Cells.Item(35, 2).Value = "兰州建设投资（控股）集团有限公司"
$ws.Cells.Item(36, 2).Value = "凉山发展（集团）集团有限公司"
$ws.Cells.Item(37, 2).Value = "临沂市建设投资集团有限公司"
$ws.Cells.Item(38, 2).Value = "牡丹江市投资集团有限公司"
$ws.Cells.Item(39, 2).Value = "青岛城市建设投资（集团）有限公司"
$ws.Cells.Item(40, 2).Value = "Shanghai Lin刚economic development (group) co., Ltd"
$ws.Cells.Item(41, 2).Value = "Shanghai Lin刚economic development (group) co., Ltd"
$ws.Cells.Item(42, 2).Value = "四川铁路投资集团有限公司"
$ws.Cells.Item(43, 2).Value = "SU纤economic development corporation"
$ws.Cells.Item(44, 2).Value = "天津滨海新区建设投资集团有限公司"
$ws.Cells.Item(45, 2).Value = "潍坊城市建设发展投资集团有限公司"
$ws.Cells.Item(46, 2).Value = "武汉贸易集团有限公司"
$ws.Cells.Item(47, 2).Value = "无锡市建设发展投资有限公司"
$ws.Cells.Item(48, 2).Value = "徐州经济技术开发区国有资产管理有限公司"
$ws.Cells.Item(49, 2).Value = "Y安泰Guofeng investment"
$ws.Cells.Item(50, 2).Value = "岳阳建设投资集团有限公司"
$ws.Cells.Item(51, 2).Value = "漳州交通发展集团有限公司"
$ws.Cells.Item(52, 2).Value = "中原资产管理有限公司"
$ws.Cells.Item(53, 2).Value = "Z宏远Y U字investment holding group co., Ltd."
$ws.Cells.Item(54, 2).Value = "Z胡海hu A法group co., Ltd."
$ws.Cells.Item(55, 2).Value = "Z胡海hu A法group co., Ltd."
$ws.Cells.Item(56, 2).Value = "诸暨市国有资产管理有限公司"
